$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# --- Row 2: OrderId (Q2) 51540586 -> 51540640 ---
# Preserve the existing look (thin top+bottom border, white solid fill) by
# stashing a copy of the current formatting in a scratch cell, changing the
# cell to Text so the numeric-looking id is stored as text (matches the
# existing column convention), writing the new id, then restoring the
# formatting from the scratch cell.
$ws.Range("Q2").Copy()
$ws.Range("AZ100").PasteSpecial(-4122)
$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = "51540640"
$ws.Range("AZ100").Copy()
$ws.Range("Q2").PasteSpecial(-4122)
$ws.Range("AZ100").Clear()

# --- Row 3: OrderId (Q3) 51540116 -> 51540772 ---
$ws.Range("Q3").Copy()
$ws.Range("AZ101").PasteSpecial(-4122)
$ws.Range("Q3").NumberFormat = "@"
$ws.Range("Q3").Value = "51540772"
$ws.Range("AZ101").Copy()
$ws.Range("Q3").PasteSpecial(-4122)
$ws.Range("AZ101").Clear()
